$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.966.81"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "3.697.10"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.42"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.90"
$ws.Range("E6").Value = "  +18.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "669.16"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  +6.90%  "
$ws.Range("E9").Value = "  +5.94%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "3.693.23"
$ws.Range("E11").Value = "  +3.53%  "
$ws.Range("E12").Value = "  +5.48%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("D15").Value = "4.386.27"
$ws.Range("E15").Value = "  +3.71%  "
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("D17").Value = "96.737.89"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.06"
$ws.Range("E18").Value = "  +17.00%  "
$ws.Range("D19").Value = "3.686.39"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.59"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.539"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.47"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "516.23"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("E25").Value = "  +5.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.97"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.71"
$ws.Range("E27").Value = "  +6.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.15"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  +11.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.186"
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("E34").Value = "  +5.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("E37").Value = "  +4.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.80"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "613.41"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.75"
$ws.Range("E40").Value = "  +26.96%  "
$ws.Range("E41").Value = "  +8.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.964"
$ws.Range("E42").Value = "  +7.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.95"
$ws.Range("E43").Value = "  +7.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0463"
$ws.Range("E45").Value = "  +9.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.21"
$ws.Range("E46").Value = "  +9.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.432"
$ws.Range("E47").Value = "  +26.46%  "
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.63"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  +5.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.82"
$ws.Range("E51").Value = "  +4.51%  "
